$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'27.857.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "'1.879.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'333.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4734"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.24%  "
$ws.Range("D8").Value = "'0.3976"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").Value = "'48.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.08036"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "'21.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").Value = "'1.912.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").Value = "'5.965"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "'7.187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'0.00001051"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "'87.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "'0.06624"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "'17.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D22").Value = "'27.904.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").Value = "'5.507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "'11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "'2.120.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "'157.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("D28").Value = "'20.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "'2.105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").Value = "'5.625"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'122.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").Value = "'0.9786"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("D33").Value = "'0.09570"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("D34").Value = "'1.465"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").Value = "'3.636"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'5.317"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "'0.06112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").Value = "'0.02263"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "'8.220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'0.6034"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'0.1910"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("D44").Value = "'10.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").Value = "'0.5716"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "'1.247"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'12.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "'3.416"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "'0.06833"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'113.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.44%  "
